# Update "想去人数" (want-to-go count) values in both the "展览" sheet and
# the aggregated "全部类型" sheet, reflecting refreshed counts for the same
# set of events in each sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Cell -> new value updates for the "展览" sheet
$exhibitionUpdates = @{
    "F4"  = 9340
    "F8"  = 257
    "F10" = 401
    "F15" = 11952
    "F30" = 2143
    "F34" = 472
    "F37" = 14
}

foreach ($cell in $exhibitionUpdates.Keys) {
    $wsExhibition.Range($cell).Value = $exhibitionUpdates[$cell]
}

# Cell -> new value updates for the "全部类型" sheet (same events, different rows)
$allTypesUpdates = @{
    "F7"  = 9340
    "F11" = 257
    "F13" = 401
    "F16" = 11952
    "F31" = 2143
    "F35" = 472
    "F38" = 14
}

foreach ($cell in $allTypesUpdates.Keys) {
    $wsAllTypes.Range($cell).Value = $allTypesUpdates[$cell]
}
